# Actualización automática 2026-01-01 08:30:07
# Monthly rollover: "VENTA MENSUAL" shifts septiembre->octubre->noviembre->diciembre
# one column to the left (oldest month drops off, a new "enero" column appears at
# the right, still empty), and the December breakdown-by-category figures on the
# "VENTAS POR GRUPO" sheet (which fed into the now-archived December column) are
# reset to 0, along with their "x de 36" participation counters.

$wb  = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: VENTAS POR GRUPO  -> zero out the cells that made up December's total
# per category, and refresh the "x de 36" counters in the summary row (38).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("I5").Value  = 0
$ws1.Range("L5").Value  = 0
$ws1.Range("M5").Value  = 0
$ws1.Range("D11").Value = 0
$ws1.Range("I11").Value = 0
$ws1.Range("L11").Value = 0
$ws1.Range("M13").Value = 0
$ws1.Range("M21").Value = 0
$ws1.Range("I26").Value = 0
$ws1.Range("M27").Value = 0
$ws1.Range("I30").Value = 0
$ws1.Range("M30").Value = 0
$ws1.Range("M34").Value = 0

$ws1.Range("D38").Value = "0 de 36"
$ws1.Range("I38").Value = "0 de 36"
$ws1.Range("L38").Value = "0 de 36"
$ws1.Range("M38").Value = "0 de 36"

# ---------------------------------------------------------------------------
# Sheet: VENTA MENSUAL -> shift month columns left (C<-D, D<-E, E<-F) and bring
# in the new "enero" column headers; column widths were tweaked too.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Columns.Item(3).ColumnWidth = 13.166666666666666   # C -> width 14
$ws2.Columns.Item(4).ColumnWidth = 14.166666666666666   # D -> width 15
$ws2.Columns.Item(5).ColumnWidth = 14.166666666666666   # E -> width 15 (unchanged)
$ws2.Columns.Item(6).ColumnWidth = 10.166666666666666   # F -> width 11

$ws2.Range("C1").Value = "octubre"
$ws2.Range("D1").Value = "noviembre"
$ws2.Range("E1").Value = "diciembre"
$ws2.Range("F1").Value = "enero"

$ws2.Range("C2").Value = 6045.88
$ws2.Range("D2").Value = 784.0599999999999
$ws2.Range("E2").Value = 2645.91

$ws2.Range("C3").Value = 1326.66
$ws2.Range("D3").Value = 1496.52
$ws2.Range("E3").Value = 0

$ws2.Range("C4").Value = 0
$ws2.Range("E4").Value = -199.54

$ws2.Range("C5").Value = 476.59
$ws2.Range("D5").Value = 286
$ws2.Range("E5").Value = 2358.37
$ws2.Range("F5").Value = 0

$ws2.Range("C9").Value = 2601.5
$ws2.Range("D9").Value = 3401.69
$ws2.Range("E9").Value = 3353.25

$ws2.Range("C11").Value = -309.07
$ws2.Range("D11").Value = 0
$ws2.Range("E11").Value = 2565.93
$ws2.Range("F11").Value = 0

$ws2.Range("C13").Value = 1314.1
$ws2.Range("D13").Value = 66.68000000000001
$ws2.Range("E13").Value = 351.48
$ws2.Range("F13").Value = 0

$ws2.Range("C16").Value = 350.18
$ws2.Range("D16").Value = -5874.77
$ws2.Range("E16").Value = 2807.2

$ws2.Range("C20").Value = 0
$ws2.Range("E20").Value = 937.86

$ws2.Range("C21").Value = 12049.42
$ws2.Range("D21").Value = 3779.22
$ws2.Range("E21").Value = 7518.26
$ws2.Range("F21").Value = 0

$ws2.Range("C22").Value = 5015.36
$ws2.Range("D22").Value = 0
$ws2.Range("E22").Value = 183.62

$ws2.Range("C24").Value = 5179.53
$ws2.Range("D24").Value = 366.34
$ws2.Range("E24").Value = 44.06

$ws2.Range("E26").Value = 44.1
$ws2.Range("F26").Value = 0

$ws2.Range("C27").Value = 6777.81
$ws2.Range("D27").Value = 0
$ws2.Range("E27").Value = 73.51000000000001
$ws2.Range("F27").Value = 0

$ws2.Range("C30").Value = 259.58
$ws2.Range("D30").Value = 998.1
$ws2.Range("E30").Value = 567.67
$ws2.Range("F30").Value = 0

$ws2.Range("C31").Value = 0

$ws2.Range("C33").Value = 2536.39
$ws2.Range("D33").Value = 557.5599999999999
$ws2.Range("E33").Value = 0

$ws2.Range("D34").Value = 59.02
$ws2.Range("E34").Value = 1676.04
$ws2.Range("F34").Value = 0

$ws2.Range("D35").Value = 2350.86
$ws2.Range("E35").Value = 0

$ws2.Range("C37").Value = 1758.38
$ws2.Range("D37").Value = 0
$ws2.Range("E37").Value = 220.32

$ws2.Range("C38").Value = 45382.31
$ws2.Range("D38").Value = 8271.279999999999
$ws2.Range("E38").Value = 25148.04
$ws2.Range("F38").Value = 0
